# Fruta / hortaliza, semanal
# Inserts a new weekly report block (3 rows) for "Choclo" / "Lluteño" priced
# on 2022-11-25 (serial 44890) ahead of the existing rows for this
# producer/product combination, pushing the rest of the table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right above the current row 807, shifting the old
# rows 807:824 down to 810:827 (dimension grows from A1:R824 to A1:R827).
$ws.Rows("807:809").Insert()

# --- Row 807: Lluteño / Primera ---------------------------------------
$ws.Range("A807").Value = 1
$ws.Range("B807").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C807").Value = "Arica y Parinacota"
$ws.Range("D807").Value = 44890
$ws.Range("E807").Value = 15
$ws.Range("F807").Value = 100112024
$ws.Range("G807").Value = "Choclo"
$ws.Range("H807").Value = "Lluteño"
$ws.Range("I807").Value = "Primera"
$ws.Range("J807").Value = 90
$ws.Range("K807").Value = 18000
$ws.Range("L807").Value = 20000
$ws.Range("M807").Value = 18667
$ws.Range("N807").Value = "`$/saco 50 unidades"
$ws.Range("O807").Value = "Región de Arica y Parinacota"
$ws.Range("P807").Value = 373
$ws.Range("Q807").Value = 50
$ws.Range("R807").Value = "Hortaliza"

# --- Row 808: Lluteño / Segunda -----------------------------------------
$ws.Range("A808").Value = 1
$ws.Range("B808").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C808").Value = "Arica y Parinacota"
$ws.Range("D808").Value = 44890
$ws.Range("E808").Value = 15
$ws.Range("F808").Value = 100112024
$ws.Range("G808").Value = "Choclo"
$ws.Range("H808").Value = "Lluteño"
$ws.Range("I808").Value = "Segunda"
$ws.Range("J808").Value = 100
$ws.Range("K808").Value = 14000
$ws.Range("L808").Value = 15000
$ws.Range("M808").Value = 14300
$ws.Range("N808").Value = "`$/saco 75 unidades"
$ws.Range("O808").Value = "Región de Arica y Parinacota"
$ws.Range("P808").Value = 191
$ws.Range("Q808").Value = 75
$ws.Range("R808").Value = "Hortaliza"

# --- Row 809: Lluteño / Tercera -----------------------------------------
$ws.Range("A809").Value = 1
$ws.Range("B809").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C809").Value = "Arica y Parinacota"
$ws.Range("D809").Value = 44890
$ws.Range("E809").Value = 15
$ws.Range("F809").Value = 100112024
$ws.Range("G809").Value = "Choclo"
$ws.Range("H809").Value = "Lluteño"
$ws.Range("I809").Value = "Tercera"
$ws.Range("J809").Value = 140
$ws.Range("K809").Value = 9000
$ws.Range("L809").Value = 10000
$ws.Range("M809").Value = 9429
$ws.Range("N809").Value = "`$/saco 100 unidades"
$ws.Range("O809").Value = "Región de Arica y Parinacota"
$ws.Range("P809").Value = 94
$ws.Range("Q809").Value = 100
$ws.Range("R809").Value = "Hortaliza"
